# The author switched the active sheet from "Pieces" to "Board indices"
# and changed which cell was selected on each sheet (visible in the
# workbookView/sheetView bits of the OOXML: activeTab, tabSelected and
# the per-sheet <selection> element).
#
# Order matters: selecting a range on a worksheet implicitly activates
# that worksheet (as in real Excel), so we set the "Pieces" selection
# first, then activate "Board indices" and select there last so it ends
# up being the final active tab.

$wb = $excel.ActiveWorkbook

$piecesSheet = $wb.Worksheets.Item("Pieces")
$boardSheet  = $wb.Worksheets.Item("Board indices")

# "Pieces" keeps a selection on I8 but is no longer the active tab.
$piecesSheet.Range("I8").Select()

# "Board indices" becomes the active sheet/tab, selection moves to B9.
$boardSheet.Activate()
$boardSheet.Range("B9").Select()
